# TechnoterBurnDownChart update - 3 Juni 2014
# - Rename the three "Pengumpulan ... requirements" tasks to mention the
#   accompanying slide deck.
# - Record 2 uninterrupted-effort hours against each of those three tasks
#   for the O-column (the "Actual" burn-down recalculates off these cells).
# - Widen column C so the longer task names are readable.
# - Leave the active selection on O21, matching where the author's cursor
#   ended up after editing the Actual burn-down total.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the three requirements-gathering tasks so they also note that a
# slide was produced for each one.
$ws.Range("C5").Value = "Pengumpulan user requirements dan pembuatan slide"
$ws.Range("C6").Value = "Pengumpulan system requirements dan pembuatan slide"
$ws.Range("C7").Value = "Pengumpulan functional requirements dan pembuatan slide"

# Log 2 hours of uninterrupted effort against each task in column O.
$ws.Range("O5").Value = 2
$ws.Range("O6").Value = 2
$ws.Range("O7").Value = 2

# Column C needs to be noticeably wider to fit the new, longer task names.
$ws.Columns("C").ColumnWidth = 58.3

# Leave the selection where the author left it after the edit.
$ws.Range("O21").Select()
